$d = $word.ActiveDocument

# Locate the bullet paragraph that ends with "...based on her request" -
# the new bullet about Jack's figures belongs directly after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*based on her request*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'based on her request' paragraph"
}

$nextPara = $target.Next()
$insPos = $nextPara.Range.Start

# Insert the new list-item paragraph as raw OOXML right before the
# following paragraph ("Checked everyone was happy..."). A trailing
# empty <w:p/> is required so the inserted paragraph keeps its own
# paragraph mark instead of merging into the following paragraph; it
# is removed again immediately afterwards.
$insertionPoint = $d.Range($insPos, $insPos)

$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Jack went through some figures he created and offered to make another for the final set of results for the report</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xmlFragment)

# Remove the now-empty placeholder paragraph that was needed to keep the
# new paragraph's mark from merging into the next paragraph.
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Jack went through some figures*") {
        $newPara = $p
        break
    }
}

if ($newPara -eq $null) {
    throw "Could not find the newly inserted paragraph"
}

$placeholder = $newPara.Next()
$placeholder.Range.Delete()
